# Natmi following Dr Hou advice
# Rebuild the Sema5a-Met sending/target cluster table as a full 3x3 cartesian
# grid over clusters {ECs, FAPs, sCs} (previously only 4 of the 9 combinations
# were present), with refreshed edge-expression statistics for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Sema5a", "Met", "ECs", 2.0, 0.6666666666666666, 0.3843716666666667, 1.153115, 0.009283037010184481, 0.009283037010184483, 3.0, 1.0, 2.254050666666667, 6.762152, 0.5759961839619929, 0.575996183961993, 0.8663932114977779, 7.797538903480001, 0.005346993893444209, 0.005346993893444211),
    @("ECs", "Sema5a", "Met", "FAPs", 2.0, 0.6666666666666666, 0.3843716666666667, 1.153115, 0.009283037010184481, 0.009283037010184483, 1.0, 0.3333333333333333, 0.05371366666666667, 0.161141, 0.01372589688605336, 0.01372589688605336, 0.02064601157944445, 0.185814104215, 0.0001274180087912093, 0.0001274180087912093),
    @("ECs", "Sema5a", "Met", "sCs", 2.0, 0.6666666666666666, 0.3843716666666667, 1.153115, 0.009283037010184481, 0.009283037010184483, 3.0, 1.0, 1.605544, 4.816632, 0.4102779191519537, 0.4102779191519537, 0.6171256231866667, 5.55413060868, 0.003808625107949062, 0.003808625107949063),
    @("FAPs", "Sema5a", "Met", "ECs", 3.0, 1.0, 26.63881666666667, 79.91645, 0.6433593900630531, 0.6433593900630532, 3.0, 1.0, 2.254050666666667, 6.762152, 0.5759961839619929, 0.575996183961993, 60.04524246671111, 540.4071822004, 0.3705725535924339, 0.370572553592434),
    @("FAPs", "Sema5a", "Met", "FAPs", 3.0, 1.0, 26.63881666666667, 79.91645, 0.6433593900630531, 0.6433593900630532, 1.0, 0.3333333333333333, 0.05371366666666667, 0.161141, 0.01372589688605336, 0.01372589688605336, 1.430868518827778, 12.87781666945, 0.00883068464867965, 0.008830684648679652),
    @("FAPs", "Sema5a", "Met", "sCs", 3.0, 1.0, 26.63881666666667, 79.91645, 0.6433593900630531, 0.6433593900630532, 3.0, 1.0, 1.605544, 4.816632, 0.4102779191519537, 0.4102779191519537, 42.76979226626667, 384.9281303964, 0.2639561518219395, 0.2639561518219396),
    @("sCs", "Sema5a", "Met", "ECs", 3.0, 1.0, 14.38262166666667, 43.147865, 0.3473575729267623, 0.3473575729267623, 3.0, 1.0, 2.254050666666667, 6.762152, 0.5759961839619929, 0.575996183961993, 32.41915795616444, 291.77242160548, 0.2000766364761148, 0.2000766364761148),
    @("sCs", "Sema5a", "Met", "FAPs", 3.0, 1.0, 14.38262166666667, 43.147865, 0.3473575729267623, 0.3473575729267623, 1.0, 0.3333333333333333, 0.05371366666666667, 0.161141, 0.01372589688605336, 0.01372589688605336, 0.7725433459961111, 6.952890113964999, 0.0047677942285825, 0.004767794228582501),
    @("sCs", "Sema5a", "Met", "sCs", 3.0, 1.0, 14.38262166666667, 43.147865, 0.3473575729267623, 0.3473575729267623, 3.0, 1.0, 1.605544, 4.816632, 0.4102779191519537, 0.4102779191519537, 23.09193192118667, 207.82738729068, 0.142513142222065, 0.1425131422220651)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
    }
}